# Apply the "harvard case classification" update.
#
# The average_doctor / average_doctor_old header pair (BP1/BQ1) is swapped:
# what used to be "average_doctor" becomes the new "average_doctor_old",
# and a freshly computed "average_doctor" takes its place. The same shift
# happens to every other "_old" comparison column (Ada_old, Avey_old,
# Buoy_old, K health_old, WebMD_old, doctor_MA_old, doctor_NJ_old,
# doctor_TH_old, average_doctor/_old): each now holds newly computed
# average/variance/std-dev figures for the stats rows 4-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column header relabel (BP/BQ swap) ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Updated stats values ---
# row 4
$ws.Range("E4").Value = 0.369
$ws.Range("F4").Value = 0.08
$ws.Range("G4").Value = 0.283
$ws.Range("N4").Value = 0.386
$ws.Range("O4").Value = 0.059
$ws.Range("P4").Value = 0.243
$ws.Range("W4").Value = 0.226
$ws.Range("X4").Value = 0.104
$ws.Range("Y4").Value = 0.322
$ws.Range("AI4").Value = 0.206
$ws.Range("AJ4").Value = 0.065
$ws.Range("AK4").Value = 0.256
$ws.Range("AU4").Value = 0.144
$ws.Range("AV4").Value = 0.027
$ws.Range("AW4").Value = 0.164
$ws.Range("BA4").Value = 1.949
$ws.Range("BB4").Value = 0.165
$ws.Range("BC4").Value = 0.407
$ws.Range("BG4").Value = 0.729
$ws.Range("BH4").Value = 0.145
$ws.Range("BI4").Value = 0.381
$ws.Range("BM4").Value = 0.6820000000000001
$ws.Range("BN4").Value = 0.09
$ws.Range("BO4").Value = 0.3
$ws.Range("BP4").Value = 0.65
$ws.Range("BQ4").Value = 0.657
# row 5
$ws.Range("E5").Value = 0.476
$ws.Range("F5").Value = 0.103
$ws.Range("G5").Value = 0.321
$ws.Range("N5").Value = 0.752
$ws.Range("O5").Value = 0.08699999999999999
$ws.Range("P5").Value = 0.295
$ws.Range("W5").Value = 0.234
$ws.Range("X5").Value = 0.116
$ws.Range("Y5").Value = 0.34
$ws.Range("AI5").Value = 0.245
$ws.Range("AJ5").Value = 0.095
$ws.Range("AK5").Value = 0.308
$ws.Range("AU5").Value = 0.285
$ws.Range("AV5").Value = 0.09
$ws.Range("AW5").Value = 0.301
$ws.Range("BA5").Value = 1.35
$ws.Range("BB5").Value = 0.08599999999999999
$ws.Range("BC5").Value = 0.293
$ws.Range("BG5").Value = 0.395
$ws.Range("BH5").Value = 0.047
$ws.Range("BI5").Value = 0.216
$ws.Range("BM5").Value = 0.5659999999999999
$ws.Range("BN5").Value = 0.077
$ws.Range("BO5").Value = 0.278
$ws.Range("BP5").Value = 0.45
$ws.Range("BQ5").Value = 0.455
# row 6
$ws.Range("E6").Value = 0.416
$ws.Range("N6").Value = 0.51
$ws.Range("W6").Value = 0.23
$ws.Range("AI6").Value = 0.224
$ws.Range("AU6").Value = 0.191
$ws.Range("BA6").Value = 1.583
$ws.Range("BG6").Value = 0.512
$ws.Range("BM6").Value = 0.619
$ws.Range("BP6").Value = 0.528
$ws.Range("BQ6").Value = 0.534
# row 7
$ws.Range("E7").Value = 0.45
$ws.Range("N7").Value = 0.632
$ws.Range("W7").Value = 0.232
$ws.Range("AI7").Value = 0.236
$ws.Range("AU7").Value = 0.238
$ws.Range("BA7").Value = 1.433
$ws.Range("BG7").Value = 0.435
$ws.Range("BM7").Value = 0.586
$ws.Range("BP7").Value = 0.478
$ws.Range("BQ7").Value = 0.483
# row 8
$ws.Range("E8").Value = 0.5
$ws.Range("F8").Value = 0.132
$ws.Range("G8").Value = 0.363
$ws.Range("N8").Value = 0.743
$ws.Range("O8").Value = 0.075
$ws.Range("P8").Value = 0.273
$ws.Range("W8").Value = 0.227
$ws.Range("X8").Value = 0.108
$ws.Range("Y8").Value = 0.328
$ws.Range("AI8").Value = 0.224
$ws.Range("AJ8").Value = 0.094
$ws.Range("AK8").Value = 0.306
$ws.Range("AU8").Value = 0.226
$ws.Range("AV8").Value = 0.07000000000000001
$ws.Range("AW8").Value = 0.265
$ws.Range("BA8").Value = 1.722
$ws.Range("BB8").Value = 0.136
$ws.Range("BC8").Value = 0.368
$ws.Range("BG8").Value = 0.555
$ws.Range("BH8").Value = 0.106
$ws.Range("BI8").Value = 0.326
$ws.Range("BM8").Value = 0.705
$ws.Range("BN8").Value = 0.073
$ws.Range("BO8").Value = 0.27
$ws.Range("BP8").Value = 0.574
$ws.Range("BQ8").Value = 0.585
# row 9
$ws.Range("E9").Value = 0.417
$ws.Range("F9").Value = 0.243
$ws.Range("G9").Value = 0.493
$ws.Range("N9").Value = 0.604
$ws.Range("O9").Value = 0.239
$ws.Range("P9").Value = 0.489
$ws.Range("W9").Value = 0.125
$ws.Range("X9").Value = 0.109
$ws.Range("Y9").Value = 0.331
$ws.Range("AI9").Value = 0.125
$ws.Range("AJ9").Value = 0.109
$ws.Range("AK9").Value = 0.331
$ws.Range("BA9").Value = 1.667
$ws.Range("BB9").Value = 0.246
$ws.Range("BC9").Value = 0.496
$ws.Range("BG9").Value = 0.583
$ws.Range("BH9").Value = 0.243
$ws.Range("BI9").Value = 0.493
$ws.Range("BM9").Value = 0.646
$ws.Range("BN9").Value = 0.229
$ws.Range("BO9").Value = 0.478
$ws.Range("BP9").Value = 0.556
$ws.Range("BQ9").Value = 0.555
# row 10
$ws.Range("E10").Value = 0.542
$ws.Range("N10").Value = 0.8120000000000001
$ws.Range("O10").Value = 0.152
$ws.Range("P10").Value = 0.39
$ws.Range("W10").Value = 0.271
$ws.Range("X10").Value = 0.197
$ws.Range("Y10").Value = 0.444
$ws.Range("AI10").Value = 0.25
$ws.Range("AJ10").Value = 0.188
$ws.Range("AK10").Value = 0.433
$ws.Range("AU10").Value = 0.208
$ws.Range("AV10").Value = 0.165
$ws.Range("AW10").Value = 0.406
$ws.Range("BA10").Value = 1.979
$ws.Range("BG10").Value = 0.625
$ws.Range("BH10").Value = 0.234
$ws.Range("BI10").Value = 0.484
$ws.Range("BM10").Value = 0.854
$ws.Range("BN10").Value = 0.125
$ws.Range("BO10").Value = 0.353
$ws.Range("BP10").Value = 0.66
$ws.Range("BQ10").Value = 0.6889999999999999
# row 11
$ws.Range("E11").Value = 0.5620000000000001
$ws.Range("F11").Value = 0.246
$ws.Range("G11").Value = 0.496
$ws.Range("N11").Value = 0.854
$ws.Range("O11").Value = 0.125
$ws.Range("P11").Value = 0.353
$ws.Range("W11").Value = 0.271
$ws.Range("X11").Value = 0.197
$ws.Range("Y11").Value = 0.444
$ws.Range("AI11").Value = 0.25
$ws.Range("AJ11").Value = 0.188
$ws.Range("AK11").Value = 0.433
$ws.Range("AU11").Value = 0.312
$ws.Range("AV11").Value = 0.215
$ws.Range("AW11").Value = 0.464
$ws.Range("BA11").Value = 1.979
$ws.Range("BG11").Value = 0.625
$ws.Range("BH11").Value = 0.234
$ws.Range("BI11").Value = 0.484
$ws.Range("BM11").Value = 0.854
$ws.Range("BN11").Value = 0.125
$ws.Range("BO11").Value = 0.353
$ws.Range("BP11").Value = 0.66
$ws.Range("BQ11").Value = 0.6889999999999999
# row 12
$ws.Range("E12").Value = 1.444
$ws.Range("F12").Value = 0.6909999999999999
$ws.Range("G12").Value = 0.831
$ws.Range("N12").Value = 1.698
$ws.Range("O12").Value = 1.699
$ws.Range("P12").Value = 1.304
$ws.Range("W12").Value = 1.846
$ws.Range("X12").Value = 0.746
$ws.Range("Y12").Value = 0.863
$ws.Range("AI12").Value = 1.917
$ws.Range("AJ12").Value = 0.91
$ws.Range("AK12").Value = 0.954
$ws.Range("AU12").Value = 2.667
$ws.Range("AV12").Value = 1.689
$ws.Range("AW12").Value = 1.3
$ws.Range("BA12").Value = 3.592
$ws.Range("BB12").Value = 0.332
$ws.Range("BC12").Value = 0.576
$ws.Range("BG12").Value = 1.067
$ws.Range("BH12").Value = 0.062
$ws.Range("BI12").Value = 0.249
$ws.Range("BM12").Value = 1.317
$ws.Range("BN12").Value = 0.363
$ws.Range("BO12").Value = 0.602
$ws.Range("BP12").Value = 1.197
$ws.Range("BQ12").Value = 1.263
# row 13
$ws.Range("E13").Value = 1.675
$ws.Range("F13").Value = 0.711
$ws.Range("G13").Value = 0.843
$ws.Range("N13").Value = 2.373
$ws.Range("O13").Value = 1.132
$ws.Range("P13").Value = 1.064
$ws.Range("W13").Value = 1.067
$ws.Range("X13").Value = 0.173
$ws.Range("Y13").Value = 0.416
$ws.Range("AI13").Value = 1.365
$ws.Range("AJ13").Value = 0.41
$ws.Range("AK13").Value = 0.64
$ws.Range("AU13").Value = 2.449
$ws.Range("AV13").Value = 1.341
$ws.Range("AW13").Value = 1.158
$ws.Range("BA13").Value = 2.468
$ws.Range("BB13").Value = 0.322
$ws.Range("BC13").Value = 0.5679999999999999
$ws.Range("BG13").Value = 0.595
$ws.Range("BH13").Value = 0.05
$ws.Range("BI13").Value = 0.223
$ws.Range("BM13").Value = 0.967
$ws.Range("BN13").Value = 0.289
$ws.Range("BO13").Value = 0.537
$ws.Range("BP13").Value = 0.823
$ws.Range("BQ13").Value = 0.787
